# Auto-generated Excel COM-interop script to update cryptos worksheet
# Applies per-row price (D) and volume % (E) updates, and a coin identity swap for rows 40/41.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '56.854.95'
$ws.Range("E2").Value = '  -0.83%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.967.39'
$ws.Range("E3").Value = '  -1.83%  '

# Row 4
$ws.Range("E4").Value = '  +0.23%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '498.73'
$ws.Range("E5").Value = '  -3.91%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '137.45'
$ws.Range("E6").Value = '  -2.43%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  +0.14%  '

# Row 8
$ws.Range("E8").Value = '  -2.27%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '7.32'
$ws.Range("E9").Value = '  -3.61%  '

# Row 10
$ws.Range("E10").Value = '  -2.92%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.356'
$ws.Range("E11").Value = '  -1.16%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '3.475.61'
$ws.Range("E12").Value = '  -1.71%  '

# Row 13
$ws.Range("E13").Value = '  -1.71%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '25.82'
$ws.Range("E14").Value = '  -1.38%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000158'
$ws.Range("E15").Value = '  -1.49%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '56.963.83'
$ws.Range("E16").Value = '  -1.30%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '6.05'
$ws.Range("E17").Value = '  +0.37%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.965.68'
$ws.Range("E18").Value = '  -1.71%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.57'
$ws.Range("E19").Value = '  -1.12%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.81'
$ws.Range("E20").Value = '  -1.94%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '319.35'
$ws.Range("E21").Value = '  -3.87%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.998'
$ws.Range("E22").Value = '  -0.22%  '

# Row 23
$ws.Range("E23").Value = '  -0.73%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.484'
$ws.Range("E24").Value = '  -0.68%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '63.15'
$ws.Range("E25").Value = '  -1.38%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.998'
$ws.Range("E26").Value = '  -0.51%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.164'
$ws.Range("E27").Value = '  -5.11%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.0₃0887'
$ws.Range("E28").Value = '  -4.53%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.50'
$ws.Range("E29").Value = '  -3.80%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.03'
$ws.Range("E30").Value = '  -1.90%  '

# Row 31
$ws.Range("E31").Value = '  -3.63%  '

# Row 32
$ws.Range("E32").Value = '  -5.71%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '20.10'
$ws.Range("E33").Value = '  -3.32%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '154.93'
$ws.Range("E34").Value = '  -1.81%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.59'
$ws.Range("E35").Value = '  -1.21%  '

# Row 36
$ws.Range("E36").Value = '  -0.77%  '

# Row 37
$ws.Range("E37").Value = '  -4.45%  '

# Row 38
$ws.Range("E38").Value = '  -1.01%  '

# Row 39
$ws.Range("E39").Value = '  -2.21%  '

# Row 40
$ws.Range("B40").Value = 'RenzoRestakedETH'
$ws.Range("C40").Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.002.65'
$ws.Range("E40").Value = '  -1.61%  '

# Row 41
$ws.Range("B41").Value = 'OKB'
$ws.Range("C41").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '37.58'
$ws.Range("E41").Value = '  +0.58%  '

# Row 42
$ws.Range("E42").Value = '  +0.07%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.72'
$ws.Range("E43").Value = '  -0.54%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.637'
$ws.Range("E44").Value = '  -2.25%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.199.36'
$ws.Range("E45").Value = '  -4.45%  '

# Row 46
$ws.Range("E46").Value = '  -3.82%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.942'
$ws.Range("E47").Value = '  -7.10%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '5.91'
$ws.Range("E48").Value = '  +0.03%  '

# Row 49
$ws.Range("E49").Value = '  -3.42%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '19.16'
$ws.Range("E50").Value = '  -1.79%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.79'
$ws.Range("E51").Value = '  -11.01%  '

